# Add a "time_taken" column (F) with per-row timestamps to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - same style as the other header cells (bold, centered, bordered)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Timestamp values for F2:F40, copied verbatim as text strings (not dates)
$timestamps = @(
    "2021-10-05 13:41:32.402123",
    "2021-10-05 13:41:32.402136",
    "2021-10-05 13:41:32.402139",
    "2021-10-05 13:41:32.402143",
    "2021-10-05 13:41:32.402146",
    "2021-10-05 13:41:32.402149",
    "2021-10-05 13:41:32.402152",
    "2021-10-05 13:41:32.402155",
    "2021-10-05 13:41:32.402158",
    "2021-10-05 13:41:32.402162",
    "2021-10-05 13:41:32.402165",
    "2021-10-05 13:41:32.402167",
    "2021-10-05 13:41:32.402170",
    "2021-10-05 13:41:32.402173",
    "2021-10-05 13:41:32.402176",
    "2021-10-05 13:41:32.402180",
    "2021-10-05 13:41:32.402183",
    "2021-10-05 13:41:32.402187",
    "2021-10-05 13:41:32.402190",
    "2021-10-05 13:41:32.402193",
    "2021-10-05 13:41:32.402196",
    "2021-10-05 13:41:32.402199",
    "2021-10-05 13:41:32.402202",
    "2021-10-05 13:41:32.402205",
    "2021-10-05 13:41:32.402209",
    "2021-10-05 13:41:32.402212",
    "2021-10-05 13:41:32.402215",
    "2021-10-05 13:41:32.402218",
    "2021-10-05 13:41:32.402221",
    "2021-10-05 13:41:32.402224",
    "2021-10-05 13:41:32.402227",
    "2021-10-05 13:41:32.402230",
    "2021-10-05 13:41:32.402234",
    "2021-10-05 13:41:32.402237",
    "2021-10-05 13:41:32.402240",
    "2021-10-05 13:41:32.402244",
    "2021-10-05 13:41:32.402247",
    "2021-10-05 13:41:32.402250",
    "2021-10-05 13:41:32.402253"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
